$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Enterprises density (per 1000 people)
$ws.Range("B11").Value = "'10.26"
$ws.Range("C11").Value = "'0.81"
$ws.Range("D11").Value = "'11.07"

# Row 12: Employment (% of total)
$ws.Range("B12").Value = "'51.27"
$ws.Range("C12").Value = "'32.75"
$ws.Range("D12").Value = "'84.02"

# Row 14: Enterprises (% of total)
$ws.Range("B14").Value = "'92.58"
$ws.Range("C14").Value = "'7.34"
$ws.Range("D14").Value = "'99.91"
